$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# Update the "current selection" row on Hoja1 to point at the 2011 data
# file / 2020 label (mirrors the row picked on Hoja2 below).
$ws1.Range("A2").Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2011.xlsx"
$ws1.Range("B2").Value = "2020"

# Move the selection on Hoja2 to the 2011 row (A11:B11), then restore
# Hoja1 as the active sheet/tab.
$ws2.Activate()
$ws2.Range("A11:B11").Select()
$ws1.Activate()
